$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/number and week-covering date range ---
$ws.Range("A8").Value = "Volume 29   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/12/2022  Through  12/18/2022"

# --- Week to Date / 28 Day / Year to Date / 2 Year crime table (rows 14-29) ---
$ws.Range("G14").Value = 2
$ws.Range("J14").Value = 11
$ws.Range("K14").Value = -54.545454545454
$ws.Range("N14").Value = -50
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "***.*"
$ws.Range("M15").Value = -13.333333333333
$ws.Range("C16").Value = "'0"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = 42.857142857142
$ws.Range("J16").Value = 59
$ws.Range("K16").Value = 83.050847457627
$ws.Range("M16").Value = -52.212389380531
$ws.Range("N16").Value = -84.187408491947
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -20
$ws.Range("I17").Value = 208
$ws.Range("J17").Value = 209
$ws.Range("K17").Value = -0.478468899521
$ws.Range("L17").Value = 1.960784313725
$ws.Range("M17").Value = 31.645569620253
$ws.Range("N17").Value = -21.212121212121
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 11
$ws.Range("H18").Value = 22.222222222222
$ws.Range("I18").Value = 95
$ws.Range("J18").Value = 73
$ws.Range("K18").Value = 30.136986301369
$ws.Range("L18").Value = -11.214953271028
$ws.Range("M18").Value = -61.382113821138
$ws.Range("N18").Value = -84.477124183006
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 75
$ws.Range("F19").Value = 17
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = -34.615384615384
$ws.Range("I19").Value = 228
$ws.Range("J19").Value = 181
$ws.Range("K19").Value = 25.966850828729
$ws.Range("L19").Value = 8.056872037914
$ws.Range("M19").Value = -3.79746835443
$ws.Range("N19").Value = -43.283582089552
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 300
$ws.Range("I20").Value = 152
$ws.Range("J20").Value = 93
$ws.Range("K20").Value = 63.440860215053
$ws.Range("L20").Value = 60
$ws.Range("M20").Value = 38.181818181818
$ws.Range("N20").Value = -92.782526115859
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = 90.90909090909
$ws.Range("F21").Value = 66
$ws.Range("G21").Value = 68
$ws.Range("H21").Value = -2.941176470588
$ws.Range("I21").Value = 809
$ws.Range("J21").Value = 641
$ws.Range("K21").Value = 26.209048361934
$ws.Range("L21").Value = 10.821917808219
$ws.Range("M21").Value = -19.180819180819
$ws.Range("N21").Value = -80.273104120946
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 10
$ws.Range("K22").Value = 25
$ws.Range("L22").Value = 11.111111111111
$ws.Range("M22").Value = -16.666666666666
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 400
$ws.Range("F23").Value = 11
$ws.Range("H23").Value = 120
$ws.Range("I23").Value = 101
$ws.Range("J23").Value = 87
$ws.Range("K23").Value = 16.091954022988
$ws.Range("L23").Value = 18.823529411764
$ws.Range("M23").Value = 71.186440677966
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 33.333333333333
$ws.Range("F24").Value = 41
$ws.Range("H24").Value = 24.242424242424
$ws.Range("I24").Value = 539
$ws.Range("J24").Value = 424
$ws.Range("K24").Value = 27.122641509434
$ws.Range("L24").Value = 17.943107221006
$ws.Range("M24").Value = 2.666666666666
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -16.666666666666
$ws.Range("F25").Value = 28
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = 27.272727272727
$ws.Range("I25").Value = 368
$ws.Range("J25").Value = 299
$ws.Range("K25").Value = 23.076923076923
$ws.Range("L25").Value = 16.455696202531
$ws.Range("M25").Value = -17.117117117117
$ws.Range("D26").Value = "'0"
$ws.Range("E26").Value = "***.*"
$ws.Range("I26").Value = 21
$ws.Range("K26").Value = -12.5
$ws.Range("L26").Value = -40
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 26
$ws.Range("J27").Value = 29
$ws.Range("K27").Value = -10.344827586206
$ws.Range("L27").Value = -10.344827586206
$ws.Range("G28").Value = 6
$ws.Range("J28").Value = 32
$ws.Range("K28").Value = -25
$ws.Range("N28").Value = -14.285714285714
$ws.Range("D29").Value = 2
$ws.Range("G29").Value = 3
$ws.Range("J29").Value = 25
$ws.Range("K29").Value = -20
$ws.Range("N29").Value = -9.090909090909
